$d = $word.ActiveDocument

# Locate the "Version Compatibility" paragraph and find the run of text
# that needs to be trimmed down to just a period, while keeping the
# bookmarkEnd (id 2 / _Hlk93267185) adjacent to the preceding bookmarkEnd
# (id 1 / _Hlk88383689) instead of after the removed text.
$target = $d.Content
$find = $target.Find
$find.ClearFormatting()
$find.Text = ", as well as PAL, 20XX, and UnclePunch."
$found = $find.Execute()

if ($found) {
    # Insert a replacement run (just a period, same formatting) right
    # after the matched text. InsertXML appends inline content at the
    # end of the supplied range without disturbing bookmark start/end
    # ordering elsewhere in the document.
    $xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="A-OTF Folk Pro H" w:eastAsia="A-OTF Folk Pro H" w:hAnsi="A-OTF Folk Pro H"/><w:sz w:val="25"/><w:szCs w:val="25"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xmlFragment)

    # Now remove the original ", as well as PAL, 20XX, and UnclePunch."
    # run text. Because it sits between bookmarkEnd id=1 and bookmarkEnd
    # id=2, deleting it leaves those two bookmarkEnd tags adjacent - the
    # bookmarkEnd for id=2 now lands right before the new "." run instead
    # of after it, matching the intended edit.
    $target.Delete()
}
